$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Runmode column (E) - set to "Yes" for every testcase row (2-33)
$ws.Range("E2:E33").Value = "Yes"

# Update the view/selection state to match the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("E2:E33").Select()
